# Add a name/date/assignment-title header block as three new, plainly
# formatted paragraphs at the very top of the document (ahead of the
# existing first paragraph), per "first module second assignment".
$d = $word.ActiveDocument

$lines = @(
    "Gabriel Sanchez Jorgensen",
    "3/24/2024",
    "Assignment Module 1.2"
)

# Build plain <w:p><w:r><w:t>...</w:t></w:r></w:p> fragments (no pPr/rPr)
# so the new paragraphs pick up no direct formatting (unlike the existing
# first paragraph, which has a firstLine indent).
$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$xml = ""
foreach ($line in $lines) {
    $xml += "<w:p $ns><w:r><w:t>$line</w:t></w:r></w:p>"
}

$target = $d.Range(0, 0)
[void]$target.InsertXML($xml)
